$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, pushing existing rows 100-117 down to 101-118.
$ws.Rows("100").Insert()

# Populate the newly inserted row 100 with its data.
$ws.Cells.Item(100, 1).Value = 5
$ws.Cells.Item(100, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(100, 3).Value = "Maule"
$ws.Cells.Item(100, 4).Value = 44889
$ws.Cells.Item(100, 5).Value = 7
$ws.Cells.Item(100, 6).Value = 100112022
$ws.Cells.Item(100, 7).Value = "Arveja Verde"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 500
$ws.Cells.Item(100, 11).Value = 18000
$ws.Cells.Item(100, 12).Value = 18000
$ws.Cells.Item(100, 13).Value = 18000
$ws.Cells.Item(100, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(100, 15).Value = "Región del Maule"
$ws.Cells.Item(100, 16).Value = 720
$ws.Cells.Item(100, 17).Value = 25
$ws.Cells.Item(100, 18).Value = "Hortaliza"
